$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting the existing rows (and the
# previously-last row, now at 159) down by one.
$ws.Rows.Item(69).Insert()

# Populate the new row 69 with the latest price-report entry.
$ws.Cells.Item(69, 1).Value = 11
$ws.Cells.Item(69, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(69, 3).Value = "Bíobío"
$ws.Cells.Item(69, 4).Value = 44482
$ws.Cells.Item(69, 5).Value = 8
$ws.Cells.Item(69, 6).Value = 100112023
$ws.Cells.Item(69, 7).Value = "Brócoli"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 850
$ws.Cells.Item(69, 11).Value = 800
$ws.Cells.Item(69, 12).Value = 900
$ws.Cells.Item(69, 13).Value = 853
$ws.Cells.Item(69, 14).Value = "$/unidad"
$ws.Cells.Item(69, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(69, 16).Value = 853
$ws.Cells.Item(69, 17).Value = 1
$ws.Cells.Item(69, 18).Value = "Hortaliza"
